$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intensivbetten")

# --- Intensivbetten gesamt (row 2) ---
# B2 and C2 hold plain digit strings ("21908" / "21771") that must stay
# text (shared string), matching the original file's storage. A direct
# Range.Value assignment would have Excel auto-coerce these into numbers,
# so build them as text-formula results first, then paste-special the
# resulting *values* back in - this keeps the cell's number format/style
# untouched (still the default "General" style) while the stored value
# stays a text string, exactly like the source workbook.
$scratch = $ws.Range("Z1")

$scratch.Formula = "=""21908"""
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163)

$scratch.Formula = "=""21771"""
$scratch.Copy()
$ws.Range("C2").PasteSpecial(-4163)

$scratch.Clear()

$ws.Range("D2").Value = "-0,6 %"

# --- Belegung durch Patient*innen mit COVID-19 (row 3) ---
$ws.Range("B3").Value = "3678`n17 %"
$ws.Range("C3").Value = "4478`n21 %"
$ws.Range("D3").Value = "21,8 %"

# --- Freie Intensivbetten (row 4) ---
$ws.Range("B4").Value = "2776`n13 %"
$ws.Range("C4").Value = "2568`n12 %"
$ws.Range("D4").Value = "-7,5 %"
